# Automation test case written for OMS-777, OMS-3046, OMS-3470
# Adds a new "Vendor_Part_Number" column (M), renames the Qty column (B)
# values from numeric indices to descriptive scenario names, and appends
# three new bulk-order scenario rows (13-15) plus one blank styled row (16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Prime the formatting for the new rows (13-16) by copying the last
#    existing data row (12), which carries the correct cell styles
#    (borders/fills) for columns A, C, D, E, F, G, H, I, J, K.
# ---------------------------------------------------------------------
$ws.Range("A12:L12").Copy($ws.Range("A13:L13"))
$ws.Range("A12:L12").Copy($ws.Range("A14:L14"))
$ws.Range("A12:L12").Copy($ws.Range("A15:L15"))
$ws.Range("A12:L12").Copy($ws.Range("A16:L16"))

# Row 16 stays blank (style-only separator row) - drop the copied values.
$ws.Range("A16:L16").ClearContents()

# Approximate the source row height (15.6) for the four new rows.
$ws.Rows.Item(13).RowHeight = 15.6
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6

# ---------------------------------------------------------------------
# 2. New "Vendor_Part_Number" column header and values.
# ---------------------------------------------------------------------
$ws.Range("M1").Value = "Vendor_Part_Number"
$ws.Range("M13").Value = "E025SLL-H"

$ws.Range("K13").Value = "46694G"
$ws.Range("K15").Value = "19853H"

# ---------------------------------------------------------------------
# 3. Rework column B ("Qty") on the existing rows from plain numbers to
#    descriptive scenario labels.
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "Two_Null_Value"
$ws.Range("B2").Value = "One_Null_Value"
$ws.Range("B4").Value = "Three_Null_Value"

$ws.Range("B13").Value = "Sku_And_Vendor_Part_No"
$ws.Range("B14").Value = "Null_Sku_and_Vendor_Part_No"
$ws.Range("B15").Value = "Sku_And_Vendor_Part_No_Mismatch"

$ws.Range("I13").Value = "Header Comment 9"
$ws.Range("I14").Value = "Header Comment 10"
$ws.Range("I15").Value = "Header Comment 11"

$ws.Range("J13").Value = "Header 9"
$ws.Range("J14").Value = "Header 10"
$ws.Range("J15").Value = "Header 11"

$ws.Range("B5").Value = "Single_Order"
$ws.Range("B6").Value = "Multiple_Order"
$ws.Range("B7").Value = "Multiple_Order"
$ws.Range("B8").Value = "Multiple_Order"
$ws.Range("B9").Value = "Multiple_Order"
$ws.Range("B10").Value = "Multiple_Order"
$ws.Range("B11").Value = "Multiple_Order"
$ws.Range("B12").Value = "Multiple_Order"

# ---------------------------------------------------------------------
# 4. Fill in the remaining cell content for the new rows 13-15 so every
#    column matches the new "bulk order" scenarios.
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "bulk_order_upload"
$ws.Range("C13").Value = "MA70"
$ws.Range("D13").Value = "CA"
$ws.Range("E13").Value = 10066860
$ws.Range("F13").Value = "JA 0911"
$ws.Range("G13").Value = "OT"
$ws.Range("H13").Value = "S"
$ws.Range("L13").Value = 1

$ws.Range("A14").Value = "bulk_order_upload"
$ws.Range("C14").Value = "MA70"
$ws.Range("D14").Value = "CA"
$ws.Range("E14").Value = 10066860
$ws.Range("F14").Value = "JA 0911"
$ws.Range("G14").Value = "OT"
$ws.Range("H14").Value = "D"
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = "E025SLL-H"

$ws.Range("A15").Value = "bulk_order_upload"
$ws.Range("C15").Value = "MA70"
$ws.Range("D15").Value = "CA"
$ws.Range("E15").Value = 10066860
$ws.Range("F15").Value = "JA 0911"
$ws.Range("G15").Value = "OT"
$ws.Range("H15").Value = "S"
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = "E025SLL-H"

# K14 intentionally stays blank, matching the source scenario data.

# ---------------------------------------------------------------------
# 5. Column sizing - add column B width, narrow column E now that it no
#    longer needs to fit the old "Header Comment 1" style long text.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 22.88671875
$ws.Columns.Item(5).ColumnWidth = 22.77734375

# ---------------------------------------------------------------------
# 6. Restore the active selection to match the edited cell (B6).
# ---------------------------------------------------------------------
$ws.Range("B6").Select()
